# Applies corrected IFRS figures for rows 2-6 (years 2014-2018) and
# clears the erroneous estimate rows 7-9 (2019E-2021E), keeping only
# the row label columns A, B, C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2 through 6, columns D..AJ
$rowData = @{
    2 = @{
        D=1196; E=2; F=2; G=-61; H=-67; I=-67; J=-1; K=613; L=261; M=351; N=349; O=2; P=76;
        Q=-15; R=-93; S=100; T=97; U=-112; V=177; W=0.14; X=-5.63; Y=-17.67; Z=-11; AA=74.41;
        AB=352.75; AC=-460; AD=-4.28; AE=2290; AF=0.86; AG=0; AH=0; AI=0; AJ=15249884
    }
    3 = @{
        D=1125; E=21; F=21; G=15; H=6; I=6; J=0; K=591; L=225; M=366; N=364; O=2; P=76;
        Q=43; R=-10; S=-14; T=9; U=33; V=156; W=1.85; X=0.49; Y=1.54; Z=0.92; AA=61.56;
        AB=370.29; AC=36; AD=66.23999999999999; AE=2385; AF=1; AG=0; AH=0; AI=0; AJ=15249884
    }
    4 = @{
        D=853; E=-4; F=-4; G=-21; H=-25; I=-24; J=0; K=599; L=260; M=339; N=337; O=2; P=76;
        Q=-24; R=-19; S=31; T=16; U=-40; V=187; W=-0.48; X=-2.88; Y=-7; Z=-4.13; AA=76.81;
        AB=337.28; AC=-161; AD=-12.48; AE=2207; AF=0.91; AG=0; AH=0; AI=0; AJ=15249884
    }
    5 = @{
        D=708; E=-51; F=-51; G=-90; H=-87; I=-87; J=0; K=870; L=602; M=268; N=255; O=14; P=76;
        Q=-40; R=-241; S=322; T=17; U=-57; V=505; W=-7.25; X=-12.35; Y=-29.57; Z=-11.91; AA=224.05;
        AB=230.58; AC=-573; AD=-6.25; AE=1670; AF=2.14; AG=0; AH=0; AI=0; AJ=15249884
    }
    6 = @{
        D=739; E=-85; F=-85; G=-145; H=-163; I=-156; K=1112; L=348; M=764; N=747; P=185;
        Q=-93; R=-158; S=351; T=14; U=-107; V=243; W=-11.43; X=-22.06; Y=-31.09; Z=-16.46; AA=45.49;
        AB=301.48; AC=-783; AD=-6.41; AE=2020; AF=2.48; AG=0; AH=0; AI=0; AJ=36959438
    }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}

# Rows 7-9 (2019E, 2020E, 2021E) had incorrect data copied in; clear all
# figure columns (D..AJ), leaving only the label columns A, B, C intact.
foreach ($r in 7..9) {
    $ws.Range("D$r`:AJ$r").ClearContents()
}
